$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows of hex-code / significance data for the Login stored
# procedure's additional validation outcomes (alarm levels FA and FB).
$ws.Range("A12").Value = "FB"
$ws.Range("A11").Value = "FA"
$ws.Range("B11").Value = "Wrong Email or PAN or National ID"
$ws.Range("B12").Value = "Password length less than 8"

$ws.Range("B12").Select()

$wb.Save()
